# Refactoring the test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 ---
$ws.Range("A2").Value = "TC_Functional_Smoke_002"
$ws.Range("B2").Value = "Pass"
$ws.Range("C2").Value = 45217.53852612023

$ws.Range("A3").Value = "TC_Functional_Smoke_003"
$ws.Range("B3").Value = "Pass"
$ws.Range("C3").Value = 45217.54047950362

$ws.Range("A4").Value = "TC_Functional_Smoke_008"
$ws.Range("B4").Value = "Pass"
$ws.Range("C4").Value = 45217.54241007673

# --- Insert a new row at position 5, pushing the old row 5 down to row 6 ---
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "TC_Functional_Smoke_003"
$ws.Range("B5").Value = "Pass"
$ws.Range("C5").Value = 45217.54501406434

# --- Update the row that was pushed down to row 6 (was the old row 5) ---
$ws.Range("A6").Value = "TC_Functional_Smoke_008"
$ws.Range("B6").Value = "Pass"
$ws.Range("C6").Value = 45217.54800098521

# --- Append new rows 7-11 ---
$ws.Range("A7").Value = "TC_Functional_Smoke_008"
$ws.Range("B7").Value = "Pass"
$ws.Range("C7").Value = 45217.55045021318
$ws.Range("C7").NumberFormat = $ws.Range("C6").NumberFormat

$ws.Range("A8").Value = "obj.TC_Functional_Sanity_002_1()"
$ws.Range("B8").Value = "Pass"
$ws.Range("C8").Value = 45217.5552775664
$ws.Range("C8").NumberFormat = $ws.Range("C6").NumberFormat

$ws.Range("A9").Value = "obj.TC_Functional_Sanity_002_2()"
$ws.Range("B9").Value = "Pass"
$ws.Range("C9").Value = 45217.56011080815
$ws.Range("C9").NumberFormat = $ws.Range("C6").NumberFormat

$ws.Range("A10").Value = "obj.TC_Functional_Sanity_002_3()"
$ws.Range("B10").Value = "Fail"
$ws.Range("C10").Value = 45217.56485854509
$ws.Range("C10").NumberFormat = $ws.Range("C6").NumberFormat

$ws.Range("A11").Value = "TC_Functional_Smoke_32"
$ws.Range("B11").Value = "Fail"
$ws.Range("C11").Value = 45217.58430443323
$ws.Range("C11").NumberFormat = $ws.Range("C6").NumberFormat
